# DEV 5 - CLI Changes for Manager Project and Enquiry
#
# A manager replies to the enquiry on row 3: the reply text is filled in
# (a new shared string "hello hello") and the reply date/time is refreshed
# to the moment the reply was sent. The enquiry date (F3) keeps its value
# but has its date/time format re-applied (same "yyyy-MM-dd HH:mm:ss"
# pattern used elsewhere in the sheet for ENQUIRY_DATE / REPLY_DATE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPLY (column E) for row 3 - manager enters a reply.
$ws.Range("E3").Value = "hello hello"

# ENQUIRY_DATE (F3) - re-apply the existing date/time number format.
$ws.Range("F3").NumberFormat = "yyyy-MM-dd HH:mm:ss"

# REPLY_DATE (G3) - record the reply timestamp and format it the same way.
$ws.Range("G3").Value = 45769.925410208336
$ws.Range("G3").NumberFormat = "yyyy-MM-dd HH:mm:ss"
